$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add month count in column I (08/24)
$ws.Range("I2").Value = 1

# Row 3: add month counts in columns I (08/24) and J (09/24)
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1

# Row 4: new row with month counts in columns I and J
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
